$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.920.08'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.50%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.125.87'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.35%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '534.04'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.56%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.72'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.41%  '

$ws.Range('E7').Value = '  -0.10%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.507'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +11.97%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.36'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.12%  '

$ws.Range('E10').Value = '  +1.95%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.418'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.79%  '

$ws.Range('E12').Value = '  +3.22%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.659.21'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.10%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.72'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.63%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000168'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.33%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '57.993.35'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.40%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.122.86'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.34%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.18'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.07%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.86'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.99%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.14'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.68%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '375.94'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +8.24%  '

$ws.Range('E22').Value = '  +0.18%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.72'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.50%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.63'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.28%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.510'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.36%  '

$ws.Range('E26').Value = '  +0.50%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.02%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0886'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.51%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.72'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +6.08%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.15'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.83%  '

$ws.Range('E31').Value = '  +0.48%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.67'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.06%  '

$ws.Range('E33').Value = '  +4.88%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.18'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.38%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '160.80'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.82%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.20'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.24%  '

$ws.Range('E37').Value = '  +6.92%  '

$ws.Range('E38').Value = '  -0.76%  '

$ws.Range('E39').Value = '  +4.46%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0674'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.99%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.19'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.91%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.555.73'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.85%  '

$ws.Range('E44').Value = '  +1.12%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0271'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.80%  '

$ws.Range('E46').Value = '  -0.04%  '

# Row 47/48: Cosmos and ONDO swap places
$ws.Range('B47').Value = 'Cosmos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.19'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.17%  '

$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.976'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.72%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0978'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +9.71%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.98'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.65%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.748'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.78%  '

